$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-08-03 Saturday" "2024-08-04 Sunday"

Replace-Text "460÷8=57, 4" "789÷7=112, 5"
Replace-Text "402÷9=44, 6" "909÷2=454, 1"
Replace-Text "668÷9=74, 2" "191÷9=21, 2"
Replace-Text "869÷7=124, 1" "249÷3=83, 0"
Replace-Text "605÷5=121, 0" "753÷7=107, 4"

Replace-Text "975÷6=162, 3" "215÷8=26, 7"
Replace-Text "462÷4=115, 2" "331÷6=55, 1"
Replace-Text "825÷9=91, 6" "403÷5=80, 3"
Replace-Text "729÷7=104, 1" "247÷6=41, 1"
Replace-Text "321÷3=107, 0" "677÷4=169, 1"

Replace-Text "147÷3=49, 0" "258÷8=32, 2"
Replace-Text "999÷2=499, 1" "605÷3=201, 2"
Replace-Text "889÷8=111, 1" "937÷2=468, 1"
Replace-Text "718÷8=89, 6" "885÷5=177, 0"
Replace-Text "231÷9=25, 6" "542÷3=180, 2"

Replace-Text "737÷2=368, 1" "610÷3=203, 1"
Replace-Text "264÷8=33, 0" "903÷2=451, 1"
Replace-Text "672÷5=134, 2" "355÷9=39, 4"
Replace-Text "294÷5=58, 4" "627÷2=313, 1"
Replace-Text "862÷5=172, 2" "709÷8=88, 5"

Replace-Text "678÷5=135, 3" "942÷6=157, 0"
Replace-Text "639÷4=159, 3" "489÷4=122, 1"
Replace-Text "642÷3=214, 0" "778÷9=86, 4"
Replace-Text "773÷4=193, 1" "254÷8=31, 6"
Replace-Text "345÷7=49, 2" "234÷4=58, 2"

Write-Host "Done"
